$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.170.36"
$ws.Range("E2").Value = "  -3.58%  "
$ws.Range("D3").Value = "1.611.59"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9984"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.58"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3778"
$ws.Range("E7").Value = "  -3.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3675"
$ws.Range("E8").Value = "  -4.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.07"
$ws.Range("E9").Value = "  -4.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9973"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.280"
$ws.Range("E11").Value = "  -6.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08106"
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.16"
$ws.Range("E13").Value = "  -4.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.643"
$ws.Range("E14").Value = "  -6.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.662"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001275"
$ws.Range("E16").Value = "  -3.53%  "
$ws.Range("D17").Value = "1.611.56"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.60"
$ws.Range("E18").Value = "  -3.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06795"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.46"
$ws.Range("E20").Value = "  -7.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.605"
$ws.Range("E21").Value = "  -4.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9978"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.08"
$ws.Range("E23").Value = "  -4.30%  "
$ws.Range("D24").Value = "23.240.36"
$ws.Range("E24").Value = "  -3.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.360"
$ws.Range("E25").Value = "  -5.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.925"
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.14"
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.81"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.277"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.61"
$ws.Range("E30").Value = "  -5.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.425"
$ws.Range("E31").Value = "  -2.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.049"
$ws.Range("D33").Value = "1.786.02"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9919"
$ws.Range("E34").Value = "  -4.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07763"
$ws.Range("E35").Value = "  -4.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02799"
$ws.Range("E36").Value = "  -6.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.344"
$ws.Range("E37").Value = "  -6.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2558"
$ws.Range("E38").Value = "  -4.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.15"
$ws.Range("E39").Value = "  -7.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08879"
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.404"
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7219"
$ws.Range("E42").Value = "  -4.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.85"
$ws.Range("E43").Value = "  -4.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.03"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6643"
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.316"
$ws.Range("E46").Value = "  -5.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9971"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.978"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08024"
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.96"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.173"
$ws.Range("E51").Value = "  -4.01%  "
